# Populate the "test" sheet (Sheet2 / ActiveSheet) with the bank/bass
# keyword-in-sentence example data used by the Cucumber-from-Excel reader.
#
# Column A holds an example sentence, column B holds the keyword that the
# sentence is being classified by ("bank" for the banking-themed sentences,
# "bass" for the fishing/music sentences). Writing column A before column B
# for each row (but writing A3/A4 before B3/B4) reproduces the exact shared
# -string insertion order of the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Tester went to bank to deposit money"
$ws.Range("B1").Value = "bank"

$ws.Range("A2").Value = "Tester went to banks of river for fishing"
$ws.Range("B2").Value = "bank"

$ws.Range("A3").Value = "I went fishing for some sea bass."
$ws.Range("A4").Value = "The bass line of the song is too weak."

$ws.Range("B3").Value = "bass"
$ws.Range("B4").Value = "bass"

# Widen column A so the longest sentence is fully visible (best-fit to the
# "Tester went to bank to deposit money" / "...banks of river..." text).
$ws.Columns.Item(1).ColumnWidth = 35.6666666666667

# Page setup: A4, portrait (matches paperSize="9" orientation="portrait").
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the cursor/selection to D13, as recorded in the saved view state.
$ws.Range("D13").Select()
